$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K"), replacing the old
# "Strike#" values with regenerated strikeout ("K") counts pulled from the
# refreshed box-score source used by the save_data regen script.
$kValues = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 1
    6 = 1
    7 = 0
    8 = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 1
    15 = 0
    16 = 2
    17 = 3
    18 = 0
    19 = 0
    20 = 0
    21 = 2
    22 = 1
    23 = 0
    24 = 2
    25 = 2
    26 = 0
    27 = 0
    28 = 2
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 1
    35 = 4
    36 = 3
    37 = 2
    38 = 0
    39 = 1
    40 = 0
    41 = 1
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 2
    47 = 1
    48 = 1
    49 = 2
    50 = 0
    51 = 2
    52 = 0
    53 = 0
    55 = 0
    56 = 2
    57 = 1
    58 = 1
    59 = 0
    60 = 1
    61 = 4
    62 = 1
    63 = 1
    64 = 3
    65 = 1
    66 = 0
    67 = 0
    68 = 0
    69 = 1
    70 = 1
    71 = 2
    72 = 2
    73 = 1
    74 = 1
    75 = 0
    76 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

